# TOLLCLASS_Designations.xlsx edit:
# "commit missing tollclass 1010001"
#
# 1) Inputs_for_tollcalib!G185:G296 (112 rows of NextGenFwy arterial segments)
#    get MAX_TOLL = 5 (previously blank).
# 2) Rows 537-547 (the last NextGenFwyR2 block) had a stale duplicate entry
#    at row 541 ("SR4 - Port Chicago..." / 401, which already exists at row 93).
#    That duplicate is dropped, rows 541-546 shift up one, the project labels
#    in column A become "NextGenFwyR2 - P8".."P18" (skipping P12), and the
#    freed last row (547) gets the missing tollclass 1010001.
# 3) The _FilterDatabase defined name / filter range grows to match the
#    sheet's real extent (K547).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs_for_tollcalib")

# --- 1) MAX_TOLL backfill, G185:G296 ---
$ws.Range("G185:G296").Value = 5

# --- 2) Fix up the tail of the NextGenFwyR2 block (rows 537-547) ---

# Column A: project label gains a "- P#" suffix for rows 537-546
$ws.Range("A537").Value = "NextGenFwyR2 - P8"
$ws.Range("A538").Value = "NextGenFwyR2 - P9"
$ws.Range("A539").Value = "NextGenFwyR2 - P10"
$ws.Range("A540").Value = "NextGenFwyR2 - P11"
$ws.Range("A541").Value = "NextGenFwyR2 - P13"
$ws.Range("A542").Value = "NextGenFwyR2 - P14"
$ws.Range("A543").Value = "NextGenFwyR2 - P15"
$ws.Range("A544").Value = "NextGenFwyR2 - P16"
$ws.Range("A545").Value = "NextGenFwyR2 - P17"
$ws.Range("A546").Value = "NextGenFwyR2 - P18"
# A547 stays "NextGenFwyR2"

# Columns B (facility_name) & C (tollclass): drop the stale row-541
# duplicate by shifting rows 542-547's old values up into 541-546, then
# put the real missing tollclass (1010001) on the newly freed row 547.
$ws.Range("B541").Value = "NA"
$ws.Range("C541").Value = 800404

$ws.Range("B542").Value = "I-580 in Contra Costa West"
$ws.Range("C542").Value = 5800504

$ws.Range("B543").Value = "NA"
$ws.Range("C543").Value = 1600501

$ws.Range("B544").Value = "SR92 in San Mateo East"
$ws.Range("C544").Value = 920202

$ws.Range("B545").Value = "SR92 in San Mateo "
$ws.Range("C545").Value = 920204

$ws.Range("B546").Value = "NA"
$ws.Range("C546").Value = 201

$ws.Range("B547").Value = "RouteNum*10000 + CountyNum*100 + SegmentNum"
$ws.Range("C547").Value = 1010001

# --- 3) Defined name / filter range update ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Inputs_for_tollcalib!_FilterDatabase") {
        $n.RefersTo = "=Inputs_for_tollcalib!`$A`$1:`$K`$547"
    }
}

# --- View state: re-freeze header row & scroll/select like the saved file ---
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 525
[void]$ws.Range("A541:XFD541").Select()

Write-Host "edit complete"
